$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) ---
$ws.Range("A1").Value = "Tuotenumero"
$ws.Range("B1").Value = "Tuotteen materiaali"
$ws.Range("C1").Value = "Poista"

# --- Row 2 ---
$ws.Range("A2").Value = "hammer123"
$ws.Range("B2").Value = "Punainen"

# --- Insert a new row at position 3, pushing the old row 3 down to row 4 ---
$ws.Rows(3).Insert() | Out-Null

# --- Fill in the newly inserted row 3 ---
$ws.Range("A3").Value = "ski1"
$ws.Range("B3").Value = "Sininen"

# --- Update what is now row 4 (previously row 3) ---
$ws.Range("A4").Value = "hammer123"
$ws.Range("B4").Value = "Aluminium"
$ws.Range("C4").Value = "X"

# --- Column widths (characters) ---
$ws.Columns("A").ColumnWidth = 11.166666666666666
$ws.Columns("B").ColumnWidth = 16.0
$ws.Columns("C").ColumnWidth = 7.333333333333333

# --- Move selection back to A1 ---
$ws.Range("A1").Select() | Out-Null

# --- Tab ratio (sheet-tab/scrollbar split) ---
try {
    $wb.Windows(1).TabRatio = 0.988
} catch {
}
